$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 45 ("Región de O'Higgins",
# 2023-10-04), pushing the previously-existing rows 45..122 down to 46..123.
$ws.Rows.Item(45).Insert()

$ws.Cells.Item(45, 1).Value = 5
$ws.Cells.Item(45, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(45, 3).Value = "Maule"
$ws.Cells.Item(45, 4).Value = 45203
$ws.Cells.Item(45, 5).Value = 7
$ws.Cells.Item(45, 6).Value = 100112026
$ws.Cells.Item(45, 7).Value = "Haba"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 300
$ws.Cells.Item(45, 11).Value = 10000
$ws.Cells.Item(45, 12).Value = 10000
$ws.Cells.Item(45, 13).Value = 10000
$ws.Cells.Item(45, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(45, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(45, 16).Value = 400
$ws.Cells.Item(45, 17).Value = 25
$ws.Cells.Item(45, 18).Value = "Hortaliza"
